# ---------------------------------------------------------------------------
# Business-agnostic template slots:
#   ミールモデル        -> 収益モデル1 (sheet2)
#   アカデミーモデル     -> 収益モデル2 (sheet3)
#   コンサルモデル       -> 収益モデル3 (sheet4)
# Plus generic-izing the row labels inside each revenue-model sheet,
# updating sample data, dropping the model-specific "LTV" row
# (replacing it with a tiny spacer row 9) and unifying the B7 formula and
# column width on the former "コンサルモデル" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1) Rename the three model sheets to generic numbered slots ------------
$wsMeal    = $wb.Worksheets.Item("ミールモデル")
$wsAcademy = $wb.Worksheets.Item("アカデミーモデル")
$wsConsul  = $wb.Worksheets.Item("コンサルモデル")

$wsMeal.Name    = "収益モデル1"
$wsAcademy.Name = "収益モデル2"
$wsConsul.Name  = "収益モデル3"

# --- 2) Fix the cross-sheet sales formulas on PL設計 (sheet1) --------------
$wsPL = $wb.Worksheets.Item("PL設計")
$wsPL.Range("B4").Formula = "='収益モデル1'!B10+'収益モデル2'!B10+'収益モデル3'!B10"
$wsPL.Range("C4").Formula = "='収益モデル1'!C10+'収益モデル2'!C10+'収益モデル3'!C10"
$wsPL.Range("D4").Formula = "='収益モデル1'!D10+'収益モデル2'!D10+'収益モデル3'!D10"
$wsPL.Range("E4").Formula = "='収益モデル1'!E10+'収益モデル2'!E10+'収益モデル3'!E10"
$wsPL.Range("F4").Formula = "='収益モデル1'!F10+'収益モデル2'!F10+'収益モデル3'!F10"

# =============================================================================
# 収益モデル1 (was ミールモデル)
# =============================================================================
$ws1 = $wb.Worksheets.Item("収益モデル1")

$ws1.Range("A1").Value = "収益モデル1（セグメント1）"

$ws1.Range("A3").Value = "顧客数/取引数"
$ws1.Range("B3").Value = 40
$ws1.Range("C3").Value = 80
$ws1.Range("D3").Value = 133
$ws1.Range("E3").Value = 213
$ws1.Range("F3").Value = 320

$ws1.Range("A4").Value = "単価（円）"
$ws1.Range("B4").Value = 4000
$ws1.Range("C4").Value = 4000
$ws1.Range("D4").Value = 4500
$ws1.Range("E4").Value = 4500
$ws1.Range("F4").Value = 5000

$ws1.Range("A5").Value = "頻度/回数（月間）"
$ws1.Range("B5").Value = 2
$ws1.Range("C5").Value = 2
$ws1.Range("D5").Value = 3
$ws1.Range("E5").Value = 3
$ws1.Range("F5").Value = 3

$ws1.Range("A6").Value = "成長率/解約率"
# B6:F6 values (0.05, 0.05, 0.04, 0.04, 0.035) are unchanged.

# Remove the model-specific "LTV" row 9 and turn it into a thin spacer row,
# like the other two model sheets.
$ws1.Rows.Item(9).Clear()
$ws1.Rows.Item(9).RowHeight = 6

# =============================================================================
# 収益モデル2 (was アカデミーモデル)
# =============================================================================
$ws2 = $wb.Worksheets.Item("収益モデル2")

$ws2.Range("A1").Value = "収益モデル2（セグメント2）"

$ws2.Range("A3").Value = "顧客数/取引数"
$ws2.Range("B3").Value = 30
$ws2.Range("C3").Value = 60
$ws2.Range("D3").Value = 100
$ws2.Range("E3").Value = 160
$ws2.Range("F3").Value = 240

$ws2.Range("A4").Value = "単価（円）"
$ws2.Range("B4").Value = 7000
$ws2.Range("C4").Value = 7000
$ws2.Range("D4").Value = 7500
$ws2.Range("E4").Value = 7500
$ws2.Range("F4").Value = 8000

$ws2.Range("A5").Value = "頻度/回数（月間）"
$ws2.Range("B5").Value = 2
$ws2.Range("C5").Value = 2
$ws2.Range("D5").Value = 3
$ws2.Range("E5").Value = 3
$ws2.Range("F5").Value = 3

$ws2.Range("A6").Value = "成長率/解約率"
$ws2.Range("B6:F6").NumberFormat = "0.0%"
$ws2.Range("B6").Value = 0.05
$ws2.Range("C6").Value = 0.05
$ws2.Range("D6").Value = 0.04
$ws2.Range("E6").Value = 0.04
$ws2.Range("F6").Value = 0.035

# Add a thin spacer row 9 (this sheet previously had no row 9 at all).
$ws2.Rows.Item(9).RowHeight = 6

# =============================================================================
# 収益モデル3 (was コンサルモデル)
# =============================================================================
$ws3 = $wb.Worksheets.Item("収益モデル3")

# Column A was a bit wider on this sheet only; bring it in line with the
# other two model sheets (width 24).
$ws3.Columns("A").ColumnWidth = $ws1.Columns("A").ColumnWidth

$ws3.Range("A1").Value = "収益モデル3（セグメント3）"

$ws3.Range("A3").Value = "顧客数/取引数"
$ws3.Range("B3").Value = 20
$ws3.Range("C3").Value = 40
$ws3.Range("D3").Value = 66
$ws3.Range("E3").Value = 106
$ws3.Range("F3").Value = 160

$ws3.Range("A4").Value = "単価（円）"
$ws3.Range("B4").Value = 10000
$ws3.Range("C4").Value = 10000
$ws3.Range("D4").Value = 10500
$ws3.Range("E4").Value = 10500
$ws3.Range("F4").Value = 11000

$ws3.Range("A5").Value = "頻度/回数（月間）"
$ws3.Range("B5").Value = 2
$ws3.Range("C5").Value = 2
$ws3.Range("D5").Value = 3
$ws3.Range("E5").Value = 3
$ws3.Range("F5").Value = 3

$ws3.Range("A6").Value = "成長率/解約率"
$ws3.Range("B6").Value = 0.05
$ws3.Range("C6").Value = 0.05
$ws3.Range("D6").Value = 0.04
$ws3.Range("E6").Value = 0.04
$ws3.Range("F6").Value = 0.035

# Monthly sales formula simplifies from a 4-factor product to a 2-factor one.
$ws3.Range("B7").Formula = "=B3*B4"
$ws3.Range("C7").Formula = "=C3*C4"
$ws3.Range("D7").Formula = "=D3*D4"
$ws3.Range("E7").Formula = "=E3*E4"
$ws3.Range("F7").Formula = "=F3*F4"

# Add a thin spacer row 9 (this sheet previously had no row 9 at all).
$ws3.Rows.Item(9).RowHeight = 6
